$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mapping of row -> new F (AUC rect) and G (AUC trap) values
# Commit message: "Changed normalization to only include first 15 seconds"
$updates = @{
    2  = @(3567.325237194857, 3570.169429357114)
    3  = @(6017.920055452805, 6022.168900184755)
    4  = @(3986.657499999912, 3989.333833333246)
    5  = @(6142.103438860858, 6147.020317420319)
    6  = @(6280.949730989804, 6284.92375836906)
    7  = @(4742.037593799409, 4745.277337315111)
    8  = @(2056.531229291038, 2058.086688204225)
    9  = @(4511.799505316266, 4513.377141169774)
    10 = @(2566.977777777925, 2570.825322165094)
    11 = @(2820.570004455493, 2822.300902205549)
    12 = @(3945.160424689995, 3948.773901369498)
    13 = @(4071.889925999785, 4076.253690634203)
    14 = @(1347.338501173534, 1347.660938766878)
    15 = @(1204.195692025733, 1204.863370760837)
    16 = @(1510.247148289058, 1511.522100760541)
    17 = @(3906.679687500216, 3909.880022321647)
    19 = @(2958.33364169675, 2960.91856367826)
    20 = @(1737.392254575702, 1739.487382140871)
    21 = @(1597.161066255186, 1599.943138619463)
    22 = @(987.8454657965844, 989.1407572007238)
    23 = @(1999.364840603998, 2001.470543437704)
    24 = @(930.7597626992998, 932.4938718947021)
    25 = @(871.442436871274, 872.7075843391132)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 6).Value = $vals[0]
    $ws.Cells.Item($row, 7).Value = $vals[1]
}
